$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.Value = '''27.418.23'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +2.05%  '

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.Value = '''1.844.92'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.66%  '

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.Value = '''1.015'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +1.29%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.Value = '''316.49'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +2.31%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +1.23%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.Value = '''0.4743'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.42%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.Value = '''0.3704'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.46%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.Value = '''0.07465'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.29%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.Value = '''0.8884'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +2.17%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.Value = '''20.53'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.58%  '

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.Value = '''1.857.87'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +5.53%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.Value = '''0.07404'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +4.69%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.Value = '''5.496'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.47%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.Value = '''93.39'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.41%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.Value = '''6.601'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.74%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.Value = '''1.015'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.27%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.Value = '''0.000008872'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +2.00%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.Value = '''14.87'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.83%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.Value = '''27.430.93'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.93%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.Value = '''5.347'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.33%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +1.66%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.Value = '''2.084.79'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +4.52%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.Value = '''1.912'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.58%  '

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.Value = '''152.55'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.86%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.Value = '''18.71'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.93%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.Value = '''2.182'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.81%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.Value = '''5.291'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.40%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.Value = '''118.13'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.13%  '

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.Value = '''0.08979'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.54%  '

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.Value = '''0.7636'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.28%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.Value = '''1.179'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.58%  '

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.Value = '''4.574'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.62%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.Value = '''2.952'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.70%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +1.33%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.Value = '''1.108'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.83%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.Value = '''0.05366'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.64%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.Value = '''0.01970'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.51%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +2.06%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.Value = '''7.344'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.31%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +2.14%  '

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.Value = '''0.5370'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +1.09%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.16%  '

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.Value = '''8.561'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.75%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.72%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.Value = '''10.57'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.13%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +1.42%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.Value = '''105.23'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.59%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.Value = '''1.685'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.06%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.Value = '''0.06341'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.99%  '
